$d = $word.ActiveDocument

# The document currently ends with a paragraph of text ("...organizational
# long-term success.") followed by a trailing empty paragraph. We need to
# insert a block of new paragraphs between those two, matching the target
# OOXML exactly (including bold heading runs, a proofErr spell-check pair,
# and bullet list runs).

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertRange = $lastPara.Range
$insertRange.Collapse(1)  # wdCollapseStart - position right before the trailing empty paragraph

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Problem Definition</w:t></w:r></w:p><w:p><w:r><w:t>The fundamental problem addressed in this project is employee attrition&#8212;i.e., the capability to detect employees who are about to leave. Employee turnover has significant far-reaching effects: direct costs such as replacing employees, indirect costs such as lost productivity, and non-monetary effects such as reduced morale among the team. Being aware of the causes and having the capability to predict attrition allows organizations to proactively respond to mitigate these effects.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The predictive challenge is formulated as a binary classification problem with the target variable, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LeaveOrNot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, being 1 for the employees who left and 0 for the employees who stayed. The solution is to create and train machine learning models on a data set that has numerous attributes (e.g., education level, age, city, experience level, gender, and whether the employee has been benched or not). The focus is not so much on accuracy as on interpretability and business relevance of the models.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Scope and Methodology</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Scope</w:t></w:r></w:p><w:p><w:r><w:t>This project centres around a publicly available dataset of anonymized employees. The information includes variables providing insight into employees' demographic, educational, and professional attributes. The scope of analysis includes:</w:t></w:r></w:p><w:p><w:r><w:t>&#8226;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Preprocessing and data cleaning</w:t></w:r></w:p><w:p><w:r><w:t>&#8226;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Exploratory Data Analysis (EDA) with a view to exposing trends and patterns</w:t></w:r></w:p><w:p><w:r><w:t>&#8226;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Feature engineering and feature selection</w:t></w:r></w:p><w:p><w:r><w:t>&#8226; Creation, training, and testing a batch of machine learning models</w:t></w:r></w:p><w:p><w:r><w:t>&#8226; Comparison of models to uniform measures</w:t></w:r></w:p><w:p><w:r><w:t>&#8226; Creation of insights and business recommendations</w:t></w:r></w:p><w:p><w:r><w:t>The project yields an end-to-end report, poster presentation at expert level, and a working predictive model. The provided data set, while incomplete, is a good mix and range of attributes from which to create worthwhile predictive models.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$null = $insertRange.InsertXML($xml)

Write-Host "Paragraph count now: $($d.Paragraphs.Count)"
